# Apply strike-through formatting to the "done" portions of Prompt 17 and
# Prompt 18 blocks (everything after the leading "Prompt: - " label, plus
# the whole following "Critérios de aceite" paragraph), matching the
# reviewed/struck-out styling used elsewhere in the document.

$d = $word.ActiveDocument

function Strike-PromptBlock($marker, $splitText) {
    $paras = $d.Paragraphs
    $idx = -1
    for ($i = 1; $i -le $paras.Count; $i++) {
        $t = $paras.Item($i).Range.Text
        if ($t -like "*$marker*") {
            $idx = $i
            break
        }
    }

    $p = $paras.Item($idx)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End

    # Locate the point right after "Prompt: - " within this paragraph so we
    # only strike the body text, leaving the leading label untouched.
    $searchRng = $d.Range($pStart, $pEnd)
    $searchRng.Find.Execute($splitText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $splitStart = $searchRng.Start

    $bodyRng = $d.Range($splitStart, $pEnd)
    $bodyRng.Font.StrikeThrough = 1

    # Strike the whole following "Critérios de aceite" paragraph, including
    # its paragraph mark, so the struck formatting carries over fully.
    $p2 = $paras.Item($idx + 1)
    $p2.Range.Font.StrikeThrough = 1
}

Strike-PromptBlock "No core, exporte" "No core, exporte"
Strike-PromptBlock "Adicione scripts no root" "Adicione scripts no root"
